# Update CalSim3 Data Extraction Init file - scenario block moved from row 28 to row 31,
# variables block moved from E253 to E258, demand units file version bump,
# demand list block (O391 -> O440), and inflow indices label updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Init")

# Scenario Indices / Names / Directory Indices / DSS path indices / Start & End Date indices
# lower-right cell references: row 28 -> row 31
$ws.Range("D5").Value = "A31"
$ws.Range("D6").Value = "B31"
$ws.Range("D7").Value = "C31"
$ws.Range("D8").Value = "G31"
$ws.Range("D9").Value = "H31"
$ws.Range("D10").Value = "I31"
$ws.Range("D11").Value = "J31"

# Variables List Indices lower-right cell: E253 -> E258
$ws.Range("D15").Value = "E258"

# Demand Units File: bump version
$ws.Range("B20").Value = "cs3rpt2022_all_demand_units_v20251111.xlsx"

# Demands List Indices lower-right cell: O391 -> O440
$ws.Range("D22").Value = "O440"

# Row 26 label: "Demands List Indices" -> "Inflow List Indices"
$ws.Range("A26").Value = "Inflow List Indices"

# Update the selected range shown in the sheet view
$ws.Range("D5:D11").Select()

$wb.Save()
